# Refactoring for read test data from excel:
# Store amountOfResults as plain numbers instead of "<n> items" text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 39

$ws.Range("B9").Select()
